$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76: ALC
$ws.Range("H76").Value = 3233.111
$ws.Range("I76").Value = 3233.0667
$ws.Range("K76").Value = 3233.0667
$ws.Range("M76").Value = -2918.0667

# Row 79: ALC
$ws.Range("H79").Value = 3233.111
$ws.Range("I79").Value = 3233.0667
$ws.Range("K79").Value = 3233.0667
$ws.Range("M79").Value = -2141.0667

# Row 97: ALC
$ws.Range("H97").Value = 1019.73334
$ws.Range("J97").Value = 1021.1429
$ws.Range("L97").Value = 3063.4287
$ws.Range("N97").Value = -4055.4287

# Row 138: ALC
$ws.Range("H138").Value = 4746.1763
$ws.Range("I138").Value = 4299.0713
$ws.Range("J138").Value = 6832.6665
$ws.Range("K138").Value = 12897.2139
$ws.Range("L138").Value = 20497.9995
$ws.Range("M138").Value = -7757.213899999999
$ws.Range("N138").Value = -30777.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 32: ARM
$ws.Range("H32").Value = 52665628
$ws.Range("I32").Value = 62532068
$ws.Range("J32").Value = 44629.332
$ws.Range("K32").Value = 62532068
$ws.Range("L32").Value = 44629.332
$ws.Range("M32").Value = -62531781
$ws.Range("N32").Value = -45203.332

# Row 37: ARM
$ws.Range("H37").Value = 25018.5
$ws.Range("J37").Value = 25018.5
$ws.Range("L37").Value = 25018.5
$ws.Range("N37").Value = -25564.5

# Row 74: ARM
$ws.Range("H74").Value = 27778736
$ws.Range("I74").Value = 217.9
$ws.Range("J74").Value = 62501884
$ws.Range("K74").Value = 217.9
$ws.Range("L74").Value = 62501884
$ws.Range("M74").Value = 656.1
$ws.Range("N74").Value = -62503632

# Row 77: ARM
$ws.Range("H77").Value = 27778736
$ws.Range("I77").Value = 217.9
$ws.Range("J77").Value = 62501884
$ws.Range("K77").Value = 1089.5
$ws.Range("L77").Value = 312509420
$ws.Range("M77").Value = 3278.5
$ws.Range("N77").Value = -312518156

$ws = $wb.Worksheets.Item("BSM")
# Row 70: BSM
$ws.Range("H70").Value = 80000
$ws.Range("J70").Value = 80000
$ws.Range("L70").Value = 80000
$ws.Range("N70").Value = -80586

# Row 73: BSM
$ws.Range("H73").Value = 80000
$ws.Range("J73").Value = 80000
$ws.Range("L73").Value = 80000
$ws.Range("N73").Value = -82028

# Row 86: BSM
$ws.Range("H86").Value = 1880.6
$ws.Range("I86").Value = 1895.4
$ws.Range("K86").Value = 1895.4
$ws.Range("M86").Value = -772.4000000000001

# Row 89: BSM
$ws.Range("H89").Value = 1880.6
$ws.Range("I89").Value = 1895.4
$ws.Range("K89").Value = 9477
$ws.Range("M89").Value = -3861

$ws = $wb.Worksheets.Item("CRP")
# Row 31: CRP
$ws.Range("H31").Value = 12159.921
$ws.Range("I31").Value = 1377.75
$ws.Range("J31").Value = 13428.412
$ws.Range("K31").Value = 1377.75
$ws.Range("L31").Value = 13428.412
$ws.Range("N31").Value = -14018.412
$ws.Range("M31").Value = -1082.75

# Row 34: CRP
$ws.Range("H34").Value = 12159.921
$ws.Range("I34").Value = 1377.75
$ws.Range("J34").Value = 13428.412
$ws.Range("K34").Value = 1377.75
$ws.Range("L34").Value = 13428.412
$ws.Range("N34").Value = -13832.412
$ws.Range("M34").Value = -1175.75

# Row 134: CRP
$ws.Range("H134").Value = 3734441.8
$ws.Range("I134").Value = 4169536.5
$ws.Range("K134").Value = 12508609.5
$ws.Range("M134").Value = -12506074.5

$ws = $wb.Worksheets.Item("CUL")
# Row 107: CUL
$ws.Range("H107").Value = 28572890
$ws.Range("I107").Value = 382.63635
$ws.Range("J107").Value = 41668624
$ws.Range("K107").Value = 1147.90905
$ws.Range("L107").Value = 125005872
$ws.Range("M107").Value = 772.09095
$ws.Range("N107").Value = -125009712

# Row 113: CUL
$ws.Range("H113").Value = 1016.7234
$ws.Range("I113").Value = 814
$ws.Range("J113").Value = 1267.7142
$ws.Range("K113").Value = 2442
$ws.Range("L113").Value = 3803.1426
$ws.Range("M113").Value = -272
$ws.Range("N113").Value = -8143.142599999999

# Row 122: CUL
$ws.Range("H122").Value = 3076.3726
$ws.Range("I122").Value = 520.6818
$ws.Range("J122").Value = 5015.1724
$ws.Range("K122").Value = 4686.1362
$ws.Range("L122").Value = 45136.55160000001
$ws.Range("M122").Value = -2236.1362
$ws.Range("N122").Value = -50036.55160000001

# Row 132: CUL
$ws.Range("H132").Value = 2535.3518
$ws.Range("I132").Value = 2474.1428
$ws.Range("J132").Value = 2601.2693
$ws.Range("K132").Value = 22267.2852
$ws.Range("L132").Value = 23411.4237
$ws.Range("M132").Value = -19737.2852
$ws.Range("N132").Value = -28471.4237

$ws = $wb.Worksheets.Item("GSM")
# Row 70: GSM
$ws.Range("H70").Value = 8200.893
$ws.Range("I70").Value = 9363.6
$ws.Range("J70").Value = 5294.125
$ws.Range("K70").Value = 9363.6
$ws.Range("L70").Value = 5294.125
$ws.Range("M70").Value = -9093.6
$ws.Range("N70").Value = -5834.125

# Row 73: GSM
$ws.Range("H73").Value = 8200.893
$ws.Range("I73").Value = 9363.6
$ws.Range("J73").Value = 5294.125
$ws.Range("K73").Value = 9363.6
$ws.Range("L73").Value = 5294.125
$ws.Range("M73").Value = -8427.6
$ws.Range("N73").Value = -7166.125

# Row 122: GSM
$ws.Range("H122").Value = 1782.6666
$ws.Range("I122").Value = 1533.3334
$ws.Range("J122").Value = 2156.6667
$ws.Range("K122").Value = 4600.0002
$ws.Range("L122").Value = 6470.000100000001
$ws.Range("M122").Value = -2150.0002
$ws.Range("N122").Value = -11370.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 99: WVR
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

# Row 113: WVR
$ws.Range("H113").Value = 1410
$ws.Range("I113").Value = 1475
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 4425
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -2255
$ws.Range("N113").Value = -5840

# Row 122: WVR
$ws.Range("H122").Value = 3118.0454
$ws.Range("I122").Value = 2662.5454
$ws.Range("J122").Value = 3573.5454
$ws.Range("K122").Value = 7987.6362
$ws.Range("L122").Value = 10720.6362
$ws.Range("M122").Value = -5537.6362
$ws.Range("N122").Value = -15620.6362

# Row 130: WVR
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 132: WVR
$ws.Range("H132").Value = 2316613
$ws.Range("I132").Value = 1655.1632
$ws.Range("J132").Value = 25003200
$ws.Range("K132").Value = 4965.4896
$ws.Range("L132").Value = 75009600
$ws.Range("M132").Value = -2435.4896
$ws.Range("N132").Value = -75014660

# Row 136: WVR
$ws.Range("H136").Value = 5116.0415
$ws.Range("I136").Value = 4990.227
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 14970.681
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -12420.681
$ws.Range("N136").Value = -24600
